# Capstone title slide: "Rossmann Sales Prediction" -> "Retail Sales Prediction"
# (slide 1, title placeholder "object 2", 2nd paragraph).
# The original paragraph is built from three runs:
#   [" "] [runA, sz=3600]  "Rossmann" [runB, sz=3600, err="1"]  " Sales Prediction" [runC, sz=3600]
# The target merges runB+runC into a single run "Retail Sales Prediction" that keeps
# runC's formatting (no err="1" misspelling flag). We reproduce that by deleting the
# "Rossmann" run outright (dropping its err="1" flag along with it) and then rewriting
# the adjoining " Sales Prediction" run's text in place, which preserves that run's own
# (clean) character formatting.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame2.TextRange

$fullText = $tr.Text

# Remove the word "Rossmann" entirely (this also removes its err="1" flagged-misspelling run).
$rossmannStart = $fullText.IndexOf("Rossmann") + 1
$rossmannRun = $tr.Characters($rossmannStart, "Rossmann".Length)
$rossmannRun.Delete()

# Re-read the text now that "Rossmann" is gone, then rewrite " Sales Prediction" in place
# (same run, same formatting) as "Retail Sales Prediction".
$fullText = $tr.Text
$salesStart = $fullText.IndexOf(" Sales Prediction") + 1
$salesRun = $tr.Characters($salesStart, " Sales Prediction".Length)
$salesRun.Text = "Retail Sales Prediction"
